$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "6_2013-02-17_298279696.jpg"
$ws.Range("A3").Value = "5_2013-02-17_298279713.jpg"
$ws.Range("A4").Value = "1_2013-02-17_298279756.jpg"
$ws.Range("A5").Value = "0_2023-07-28_457239382.jpg"
$ws.Range("A6").Value = "0_2023-07-28_457239377.jpg"
